# The presentation currently uses the "Integral" (Red Violet) design theme
# for its slide master (ppt/theme/theme1.xml) and the stock "Office Theme"
# colours for its notes master (ppt/theme/theme2.xml). The author switched
# the presentation's design back to the default "Office Theme" palette, so
# the slide master's theme colours need to become the standard Office
# palette (dk2/lt2/accent1-6/hlink/folHlink) instead of the Red Violet ones.
#
# PowerPoint exposes theme colours through Master.ColorScheme (a legacy
# 12-slot RGBColor indexer: 1=dk1 2=lt1 3=dk2 4=lt2 5-10=accent1-6
# 11=hlink 12=folHlink). RGB values are assigned in the standard VBA
# 0xBBGGRR packed-long form.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

$cs.Colors(1).RGB  = 0x000000   # dk1      -> 000000
$cs.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$cs.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$cs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$cs.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$cs.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$cs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$cs.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$cs.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$cs.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$cs.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$cs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
